$d = $word.ActiveDocument

# --- Paragraph 1: "An images folder ... GutHub." ---
# Locate the paragraph that currently holds the "An images folder" / "GutHub." text.
$r1 = $d.Content
$r1.Find.Execute("An images folder", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $r1.Paragraphs(1)
$range1 = $para1.Range

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p w14:paraId="1A509569" w14:textId="62009EC1" w:rsidR="00781F93" w:rsidRDefault="00781F93" w:rsidP="00781F93"><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="6"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="180" w:after="180" w:line="240" w:lineRule="auto"/><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">An images folder </w:t></w:r><w:r w:rsidR="00F03943"><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr><w:t>which contains a dummy image to allow it to be managed by GutHub</w:t></w:r><w:r><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">  You do not need to populate this folder (unless you want to)</w:t></w:r><w:r><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$range1.InsertXML($xml1)

# --- Paragraph 2: "An app.py file ..." ---
# Locate the paragraph that currently holds the "An app.py file" text.
$r2 = $d.Content
$r2.Find.Execute("An app.py file that contains the overall structure", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $r2.Paragraphs(1)
$range2 = $para2.Range

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p w14:paraId="2080A1C8" w14:textId="77777777" w:rsidR="00781F93" w:rsidRDefault="00781F93" w:rsidP="00781F93"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="180" w:after="180" w:line="240" w:lineRule="auto"/><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="2D3B45"/><w:sz w:val="24"/></w:rPr><w:t>An app.py file that contains the overall structure for your Python program</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$range2.InsertXML($xml2)

Write-Output "Edit complete"
